$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("global_settings")

# FTUEX A/B test options (UI SETTINGS block, row 11: D=showNextDragonInXpBar,
# E=showUnlockProgressionText, F=mapAsButton) should all be set to FALSE.
$ws.Range("D11").Value = $false
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = $false
